$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Feuil1")

# Delete column B (Colonne2) - the whole column, data shifts left (C->B, D->C, ...)
$ws.Range("B:B").EntireColumn.Delete()

# The bound Excel Table ("Tableau1") needs to be shrunk to match the new data
# extent (A1:F25, 6 columns instead of 7).
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:F25"))

# Realign the table's column definitions with the (now shifted) header
# row text so the table column names match the data again.
$ws.Range("B1").Value = "Colonne3"
$ws.Range("C1").Value = "Colonne4"
$ws.Range("D1").Value = "Colonne5"
$ws.Range("E1").Value = "Colonne6"
$ws.Range("F1").Value = "Colonne7"

$ws.Range("P23").Select()
